$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$c = $t.Cell(1, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "2+49="
$c = $t.Cell(1, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "4+62="
$c = $t.Cell(1, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "4+23="
$c = $t.Cell(1, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "61-46="
$c = $t.Cell(1, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "14+44="
$c = $t.Cell(2, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "96-44="
$c = $t.Cell(2, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "0+9="
$c = $t.Cell(2, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "66-10="
$c = $t.Cell(2, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "72+7="
$c = $t.Cell(2, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "24+31="
$c = $t.Cell(3, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "48-1="
$c = $t.Cell(3, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "81-15="
$c = $t.Cell(3, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "84-46="
$c = $t.Cell(3, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "86+13="
$c = $t.Cell(3, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "79-12="
$c = $t.Cell(4, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "58-36="
$c = $t.Cell(4, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "14+0="
$c = $t.Cell(4, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "85-59="
$c = $t.Cell(4, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "53+34="
$c = $t.Cell(4, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "56+0="
$c = $t.Cell(5, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "24+0="
$c = $t.Cell(5, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "73+22="
$c = $t.Cell(5, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "37-22="
$c = $t.Cell(5, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "84-61="
$c = $t.Cell(5, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "61+16="
$c = $t.Cell(6, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "79+10="
$c = $t.Cell(6, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "12+50="
$c = $t.Cell(6, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "29+59="
$c = $t.Cell(6, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "86-26="
$c = $t.Cell(6, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "98-82="
$c = $t.Cell(7, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "28+1="
$c = $t.Cell(7, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "5+11="
$c = $t.Cell(7, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "66-18="
$c = $t.Cell(7, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "24+55="
$c = $t.Cell(7, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "50-48="
$c = $t.Cell(8, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "2+7="
$c = $t.Cell(8, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "97-74="
$c = $t.Cell(8, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "88-67="
$c = $t.Cell(8, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "95-31="
$c = $t.Cell(8, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "69-5="
$c = $t.Cell(9, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "43-2="
$c = $t.Cell(9, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "82-8="
$c = $t.Cell(9, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "6+92="
$c = $t.Cell(9, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "27+29="
$c = $t.Cell(9, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "24-3="
$c = $t.Cell(10, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "45+14="
$c = $t.Cell(10, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "44+48="
$c = $t.Cell(10, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "13+6="
$c = $t.Cell(10, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "76-18="
$c = $t.Cell(10, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "15+30="
$c = $t.Cell(11, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "50-49="
$c = $t.Cell(11, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "70+2="
$c = $t.Cell(11, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "85-43="
$c = $t.Cell(11, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "14+44="
$c = $t.Cell(11, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "12+45="
$c = $t.Cell(12, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "8+78="
$c = $t.Cell(12, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "82-21="
$c = $t.Cell(12, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "3+16="
$c = $t.Cell(12, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "93-85="
$c = $t.Cell(12, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "11+19="
$c = $t.Cell(13, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "10-7="
$c = $t.Cell(13, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "3+13="
$c = $t.Cell(13, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "89-36="
$c = $t.Cell(13, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "19+36="
$c = $t.Cell(13, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "13-2="
$c = $t.Cell(14, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "37+2="
$c = $t.Cell(14, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "98-9="
$c = $t.Cell(14, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "47-36="
$c = $t.Cell(14, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "51+17="
$c = $t.Cell(14, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "83-57="
$c = $t.Cell(15, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "94-84="
$c = $t.Cell(15, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "55-39="
$c = $t.Cell(15, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "17-2="
$c = $t.Cell(15, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "87-26="
$c = $t.Cell(15, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "0+73="
$c = $t.Cell(16, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "68+19="
$c = $t.Cell(16, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "72-14="
$c = $t.Cell(16, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "66-26="
$c = $t.Cell(16, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "87-68="
$c = $t.Cell(16, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "77-42="
$c = $t.Cell(17, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "1+26="
$c = $t.Cell(17, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "64+30="
$c = $t.Cell(17, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "24+28="
$c = $t.Cell(17, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "66-27="
$c = $t.Cell(17, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "86-49="
$c = $t.Cell(18, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "63+27="
$c = $t.Cell(18, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "65-52="
$c = $t.Cell(18, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "33+22="
$c = $t.Cell(18, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "56-34="
$c = $t.Cell(18, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "10+0="
$c = $t.Cell(19, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "8+6="
$c = $t.Cell(19, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "66-36="
$c = $t.Cell(19, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "53-49="
$c = $t.Cell(19, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "0+95="
$c = $t.Cell(19, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "53-41="
$c = $t.Cell(20, 1)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "58-24="
$c = $t.Cell(20, 2)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "46-12="
$c = $t.Cell(20, 3)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "32-8="
$c = $t.Cell(20, 4)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "32+21="
$c = $t.Cell(20, 5)
$r = $c.Range
$r.End = $r.End - 1
$r.Text = "2+60="
